$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.243164333333333
$ws.Range("H2").Value = 3.729493
$ws.Range("I2").Value = 0.5291000614577227
$ws.Range("J2").Value = 0.5291000614577227
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.6187893333333333
$ws.Range("N2").Value = 1.856368
$ws.Range("O2").Value = 0.1765034355725207
$ws.Range("P2").Value = 0.1765034355725208
$ws.Range("Q2").Value = 0.769256829047111
$ws.Range("R2").Value = 6.923311461424
$ws.Range("S2").Value = 0.09338797860891992
$ws.Range("T2").Value = 0.09338797860891994
$ws.Range("G3").Value = 1.243164333333333
$ws.Range("H3").Value = 3.729493
$ws.Range("I3").Value = 0.5291000614577227
$ws.Range("J3").Value = 0.5291000614577227
$ws.Range("N3").Value = 6.701951000000001
$ws.Range("O3").Value = 0.6372213788099619
$ws.Range("P3").Value = 0.6372213788099619
$ws.Range("Q3").Value = 2.777208815649223
$ws.Range("R3").Value = 24.99487934084301
$ws.Range("S3").Value = 0.3371538706905256
$ws.Range("T3").Value = 0.3371538706905256
$ws.Range("G4").Value = 1.243164333333333
$ws.Range("H4").Value = 3.729493
$ws.Range("I4").Value = 0.5291000614577227
$ws.Range("J4").Value = 0.5291000614577227
$ws.Range("M4").Value = 0.207158
$ws.Range("N4").Value = 0.6214740000000001
$ws.Range("O4").Value = 0.05908973658186135
$ws.Range("P4").Value = 0.05908973658186135
$ws.Range("Q4").Value = 0.2575314369646667
$ws.Range("R4").Value = 2.317782932682
$ws.Range("S4").Value = 0.03126438325698348
$ws.Range("T4").Value = 0.03126438325698348
$ws.Range("G5").Value = 1.243164333333333
$ws.Range("H5").Value = 3.729493
$ws.Range("I5").Value = 0.5291000614577227
$ws.Range("J5").Value = 0.5291000614577227
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.1889926666666667
$ws.Range("N5").Value = 0.566978
$ws.Range("O5").Value = 0.05390825789608347
$ws.Range("P5").Value = 0.05390825789608347
$ws.Range("Q5").Value = 0.2349489424615556
$ws.Range("R5").Value = 2.114540482154
$ws.Range("S5").Value = 0.02852286256589653
$ws.Range("T5").Value = 0.02852286256589653
$ws.Range("G6").Value = 1.243164333333333
$ws.Range("H6").Value = 3.729493
$ws.Range("I6").Value = 0.5291000614577227
$ws.Range("J6").Value = 0.5291000614577227
$ws.Range("M6").Value = 0.2568966666666667
$ws.Range("N6").Value = 0.77069
$ws.Range("O6").Value = 0.07327719113957255
$ws.Range("P6").Value = 0.07327719113957255
$ws.Range("Q6").Value = 0.3193647733522222
$ws.Range("R6").Value = 2.87428296017
$ws.Range("S6").Value = 0.03877096633539713
$ws.Range("T6").Value = 0.03877096633539713
$ws.Range("I7").Value = 0.3764855829716142
$ws.Range("J7").Value = 0.3764855829716142
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.6187893333333333
$ws.Range("N7").Value = 1.856368
$ws.Range("O7").Value = 0.1765034355725207
$ws.Range("P7").Value = 0.1765034355725208
$ws.Range("Q7").Value = 0.5473711436373333
$ws.Range("R7").Value = 4.926340292736
$ws.Range("S7").Value = 0.06645099883801323
$ws.Range("T7").Value = 0.06645099883801323
$ws.Range("I8").Value = 0.3764855829716142
$ws.Range("J8").Value = 0.3764855829716142
$ws.Range("N8").Value = 6.701951000000001
$ws.Range("O8").Value = 0.6372213788099619
$ws.Range("P8").Value = 0.6372213788099619
$ws.Range("S8").Value = 0.2399046622832443
$ws.Range("T8").Value = 0.2399046622832443
$ws.Range("I9").Value = 0.3764855829716142
$ws.Range("J9").Value = 0.3764855829716142
$ws.Range("M9").Value = 0.207158
$ws.Range("N9").Value = 0.6214740000000001
$ws.Range("O9").Value = 0.05908973658186135
$ws.Range("P9").Value = 0.05908973658186135
$ws.Range("Q9").Value = 0.183248652272
$ws.Range("R9").Value = 1.649237870448
$ws.Range("S9").Value = 0.02224643392466119
$ws.Range("T9").Value = 0.02224643392466119
$ws.Range("I10").Value = 0.3764855829716142
$ws.Range("J10").Value = 0.3764855829716142
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.1889926666666667
$ws.Range("N10").Value = 0.566978
$ws.Range("O10").Value = 0.05390825789608347
$ws.Range("P10").Value = 0.05390825789608347
$ws.Range("Q10").Value = 0.1671798890506667
$ws.Range("R10").Value = 1.504619001456
$ws.Range("S10").Value = 0.02029568190099111
$ws.Range("T10").Value = 0.02029568190099111
$ws.Range("I11").Value = 0.3764855829716142
$ws.Range("J11").Value = 0.3764855829716142
$ws.Range("M11").Value = 0.2568966666666667
$ws.Range("N11").Value = 0.77069
$ws.Range("O11").Value = 0.07327719113957255
$ws.Range("P11").Value = 0.07327719113957255
$ws.Range("Q11").Value = 0.2272466809866666
$ws.Range("R11").Value = 2.04522012888
$ws.Range("S11").Value = 0.02758780602470438
$ws.Range("T11").Value = 0.02758780602470438
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.2218343333333333
$ws.Range("H12").Value = 0.665503
$ws.Range("I12").Value = 0.09441435557066305
$ws.Range("J12").Value = 0.09441435557066304
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.6187893333333333
$ws.Range("N12").Value = 1.856368
$ws.Range("O12").Value = 0.1765034355725207
$ws.Range("P12").Value = 0.1765034355725208
$ws.Range("Q12").Value = 0.1372687192337778
$ws.Range("R12").Value = 1.235418473104
$ws.Range("S12").Value = 0.01666445812558759
$ws.Range("T12").Value = 0.01666445812558759
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.2218343333333333
$ws.Range("H13").Value = 0.665503
$ws.Range("I13").Value = 0.09441435557066305
$ws.Range("J13").Value = 0.09441435557066304
$ws.Range("N13").Value = 6.701951000000001
$ws.Range("O13").Value = 0.6372213788099619
$ws.Range("P13").Value = 0.6372213788099619
$ws.Range("Q13").Value = 0.4955742773725556
$ws.Range("R13").Value = 4.460168496353001
$ws.Range("S13").Value = 0.06016284583619191
$ws.Range("T13").Value = 0.06016284583619191
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.2218343333333333
$ws.Range("H14").Value = 0.665503
$ws.Range("I14").Value = 0.09441435557066305
$ws.Range("J14").Value = 0.09441435557066304
$ws.Range("M14").Value = 0.207158
$ws.Range("N14").Value = 0.6214740000000001
$ws.Range("O14").Value = 0.05908973658186135
$ws.Range("P14").Value = 0.05908973658186135
$ws.Range("Q14").Value = 0.04595475682466667
$ws.Range("R14").Value = 0.4135928114220001
$ws.Range("S14").Value = 0.005578919400216673
$ws.Range("T14").Value = 0.005578919400216672
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.2218343333333333
$ws.Range("H15").Value = 0.665503
$ws.Range("I15").Value = 0.09441435557066305
$ws.Range("J15").Value = 0.09441435557066304
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 0.1889926666666667
$ws.Range("N15").Value = 0.566978
$ws.Range("O15").Value = 0.05390825789608347
$ws.Range("P15").Value = 0.05390825789608347
$ws.Range("Q15").Value = 0.04192506221488889
$ws.Range("R15").Value = 0.377325559934
$ws.Range("S15").Value = 0.005089713429195829
$ws.Range("T15").Value = 0.005089713429195827
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.2218343333333333
$ws.Range("H16").Value = 0.665503
$ws.Range("I16").Value = 0.09441435557066305
$ws.Range("J16").Value = 0.09441435557066304
$ws.Range("M16").Value = 0.2568966666666667
$ws.Range("N16").Value = 0.77069
$ws.Range("O16").Value = 0.07327719113957255
$ws.Range("P16").Value = 0.07327719113957255
$ws.Range("Q16").Value = 0.05698850078555556
$ws.Range("R16").Value = 0.51289650707
$ws.Range("S16").Value = 0.006918418779471043
$ws.Range("T16").Value = 0.006918418779471042
